# Add ten new Mac-Address rows (147-156) to the master-reg_center_device_h sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 147
$startDeviceId = 3000166
$count = 10

for ($i = 0; $i -lt $count; $i++) {
    $row = $startRow + $i
    $deviceId = $startDeviceId + $i

    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = $deviceId
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
    $ws.Cells.Item($row, 7).Value = "now()"
}

# Match the saved view state: scrolled so row 140 is at the top, with the
# last newly-added cell (E155) selected.
$ws.Application.ActiveWindow.ScrollRow = 140
$ws.Range("E155").Select()
